$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows (Indice, Distancia, max, min, Tempo)
$data = @(
    @(0, 1266, 1266, 1266, 0.01073165734608968),
    @(1, 911, 911, 911, 0.01048827966054281),
    @(2, 1233, 1233, 1233, 0.01220569610595703),
    @(3, 1430, 1430, 1430, 0.01220326423645019),
    @(4, 1198, 1198, 1198, 0.01245652834574382),
    @(5, 1284, 1284, 1284, 0.01223506132761637),
    @(6, 1545, 1545, 1545, 0.01256210803985596),
    @(7, 1733, 1733, 1733, 0.01212185223897298),
    @(8, 1770, 1770, 1770, 0.01214772860209147),
    @(9, 1741, 1741, 1741, 0.0121092955271403)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex++
}
